# Update the timestamp strings recorded by the handback-status report.
# These cells hold plain text (not real dates) that just happen to look
# like "yyyy-mm-dd HH:mm:ss" timestamps, so we update the text directly.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 07:05:19"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 07:05:14"
$wsZhCn.Range("K2").Value = "2016-08-18 07:05:31"

# de-de sheet: "Correspond Handoff Datetime" (H2, shares the same
# timestamp text as Overview!G2) and "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-18 07:05:19"
$wsDeDe.Range("K2").Value = "2016-08-18 07:05:39"
